$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.962739884853363
$ws.Range("B1").Value = 0.8480185866355896
$ws.Range("C1").Value = 3.721126079559326
$ws.Range("D1").Value = 2.844532012939453
$ws.Range("E1").Value = 1.29451322555542
